# Auto-generated: apply scheduled-runner price/profit updates to the Leve tables.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 248.11905  # H55: 251.78049 -> 248.11905
$ws.Cells.Item(55, 9).Value = 249.58974  # I55: 253.52632 -> 249.58974
$ws.Cells.Item(55, 10).Value = 229  # J55: 229.66667 -> 229
$ws.Cells.Item(55, 11).Value = 249.58974  # K55: 253.52632 -> 249.58974
$ws.Cells.Item(55, 12).Value = 229  # L55: 229.66667 -> 229
$ws.Cells.Item(55, 13).Value = -35.58974000000001  # M55: -39.52632 -> -35.58974000000001
$ws.Cells.Item(55, 14).Value = -657  # N55: -657.6666700000001 -> -657
$ws.Cells.Item(92, 8).Value = 459.58334  # H92: 430.30768 -> 459.58334
$ws.Cells.Item(92, 9).Value = 500.9091  # I92: 465.75 -> 500.9091
$ws.Cells.Item(92, 11).Value = 500.9091  # K92: 465.75 -> 500.9091
$ws.Cells.Item(92, 13).Value = 747.0908999999999  # M92: 782.25 -> 747.0908999999999
$ws.Cells.Item(99, 8).Value = 83347460  # H99: 90924460 -> 83347460
$ws.Cells.Item(99, 10).Value = 500000200  # J99: 1000000000 -> 500000200
$ws.Cells.Item(99, 12).Value = 1500000600  # L99: 3000000000 -> 1500000600
$ws.Cells.Item(99, 14).Value = -1500003596  # N99: -3000002996 -> -1500003596
$ws.Cells.Item(116, 8).Value = 4426.143  # H116: 4426.2856 -> 4426.143
$ws.Cells.Item(116, 9).Value = 3994.3333  # I116: 3994.6667 -> 3994.3333
$ws.Cells.Item(116, 11).Value = 3994.3333  # K116: 3994.6667 -> 3994.3333
$ws.Cells.Item(116, 13).Value = -552.3332999999998  # M116: -552.6667000000002 -> -552.3332999999998
$ws.Cells.Item(132, 8).Value = 2269.4211  # H132: 2840.6 -> 2269.4211
$ws.Cells.Item(132, 9).Value = 1889.3529  # I132: 2431.4614 -> 1889.3529
$ws.Cells.Item(132, 11).Value = 5668.0587  # K132: 7294.3842 -> 5668.0587
$ws.Cells.Item(132, 13).Value = -3138.0587  # M132: -4764.3842 -> -3138.0587

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2042.5834  # H2: 2155.6365 -> 2042.5834
$ws.Cells.Item(2, 9).Value = 1683.2727  # I2: 1771.7 -> 1683.2727
$ws.Cells.Item(2, 11).Value = 1683.2727  # K2: 1771.7 -> 1683.2727
$ws.Cells.Item(2, 13).Value = -1570.2727  # M2: -1658.7 -> -1570.2727
$ws.Cells.Item(32, 8).Value = 60625.7  # H32: 66767.55499999999 -> 60625.7
$ws.Cells.Item(32, 9).Value = 10753.8125  # I32: 11424.2 -> 10753.8125
$ws.Cells.Item(32, 10).Value = 260113.25  # J32: 343484.34 -> 260113.25
$ws.Cells.Item(32, 11).Value = 10753.8125  # K32: 11424.2 -> 10753.8125
$ws.Cells.Item(32, 12).Value = 260113.25  # L32: 343484.34 -> 260113.25
$ws.Cells.Item(32, 13).Value = -10466.8125  # M32: -11137.2 -> -10466.8125
$ws.Cells.Item(32, 14).Value = -260687.25  # N32: -344058.34 -> -260687.25
$ws.Cells.Item(45, 8).Value = 6505.885  # H45: 6327.2964 -> 6505.885
$ws.Cells.Item(45, 9).Value = 9699.154  # I45: 9703 -> 9699.154
$ws.Cells.Item(45, 10).Value = 3312.6155  # J45: 3192.7144 -> 3312.6155
$ws.Cells.Item(45, 11).Value = 9699.154  # K45: 9703 -> 9699.154
$ws.Cells.Item(45, 12).Value = 3312.6155  # L45: 3192.7144 -> 3312.6155
$ws.Cells.Item(45, 13).Value = -9322.154  # M45: -9326 -> -9322.154
$ws.Cells.Item(45, 14).Value = -4066.6155  # N45: -3946.7144 -> -4066.6155
$ws.Cells.Item(61, 8).Value = 2296.9656  # H61: 2538.24 -> 2296.9656
$ws.Cells.Item(61, 9).Value = 2083.4583  # I61: 2247.8 -> 2083.4583
$ws.Cells.Item(61, 10).Value = 3321.8  # J61: 3700 -> 3321.8
$ws.Cells.Item(61, 11).Value = 2083.4583  # K61: 2247.8 -> 2083.4583
$ws.Cells.Item(61, 12).Value = 3321.8  # L61: 3700 -> 3321.8
$ws.Cells.Item(61, 13).Value = -1871.4583  # M61: -2035.8 -> -1871.4583
$ws.Cells.Item(61, 14).Value = -3745.8  # N61: -4124 -> -3745.8
$ws.Cells.Item(74, 8).Value = 1616.8914  # H74: 1644.7446 -> 1616.8914
$ws.Cells.Item(74, 9).Value = 1560.6052  # I74: 1589.8379 -> 1560.6052
$ws.Cells.Item(74, 10).Value = 1884.25  # J74: 1847.9 -> 1884.25
$ws.Cells.Item(74, 11).Value = 1560.6052  # K74: 1589.8379 -> 1560.6052
$ws.Cells.Item(74, 12).Value = 1884.25  # L74: 1847.9 -> 1884.25
$ws.Cells.Item(74, 13).Value = -686.6052  # M74: -715.8379 -> -686.6052
$ws.Cells.Item(74, 14).Value = -3632.25  # N74: -3595.9 -> -3632.25
$ws.Cells.Item(77, 8).Value = 1616.8914  # H77: 1644.7446 -> 1616.8914
$ws.Cells.Item(77, 9).Value = 1560.6052  # I77: 1589.8379 -> 1560.6052
$ws.Cells.Item(77, 10).Value = 1884.25  # J77: 1847.9 -> 1884.25
$ws.Cells.Item(77, 11).Value = 7803.026  # K77: 7949.1895 -> 7803.026
$ws.Cells.Item(77, 12).Value = 9421.25  # L77: 9239.5 -> 9421.25
$ws.Cells.Item(77, 13).Value = -3435.026  # M77: -3581.1895 -> -3435.026
$ws.Cells.Item(77, 14).Value = -18157.25  # N77: -17975.5 -> -18157.25
$ws.Cells.Item(110, 8).Value = 1936.3636  # H110: 2050 -> 1936.3636
$ws.Cells.Item(110, 9).Value = 1928.7  # I110: 2054.111 -> 1928.7
$ws.Cells.Item(110, 11).Value = 1928.7  # K110: 2054.111 -> 1928.7
$ws.Cells.Item(110, 13).Value = 116.3  # M110: -9.110999999999876 -> 116.3
$ws.Cells.Item(116, 8).Value = 2042.5834  # H116: 2155.6365 -> 2042.5834
$ws.Cells.Item(116, 9).Value = 1683.2727  # I116: 1771.7 -> 1683.2727
$ws.Cells.Item(116, 11).Value = 1683.2727  # K116: 1771.7 -> 1683.2727
$ws.Cells.Item(116, 13).Value = 610.7273  # M116: 522.3 -> 610.7273
$ws.Cells.Item(122, 8).Value = 1945.6154  # H122: 2081.838 -> 1945.6154
$ws.Cells.Item(122, 9).Value = 1837.9412  # I122: 1972.258 -> 1837.9412
$ws.Cells.Item(122, 10).Value = 2677.8  # J122: 2648 -> 2677.8
$ws.Cells.Item(122, 11).Value = 5513.8236  # K122: 5916.774 -> 5513.8236
$ws.Cells.Item(122, 12).Value = 8033.400000000001  # L122: 7944 -> 8033.400000000001
$ws.Cells.Item(122, 13).Value = -3063.8236  # M122: -3466.774 -> -3063.8236
$ws.Cells.Item(122, 14).Value = -12933.4  # N122: -12844 -> -12933.4
$ws.Cells.Item(136, 8).Value = 2296.9656  # H136: 2538.24 -> 2296.9656
$ws.Cells.Item(136, 9).Value = 2083.4583  # I136: 2247.8 -> 2083.4583
$ws.Cells.Item(136, 10).Value = 3321.8  # J136: 3700 -> 3321.8
$ws.Cells.Item(136, 11).Value = 6250.374899999999  # K136: 6743.400000000001 -> 6250.374899999999
$ws.Cells.Item(136, 12).Value = 9965.400000000001  # L136: 11100 -> 9965.400000000001
$ws.Cells.Item(136, 13).Value = -3700.374899999999  # M136: -4193.400000000001 -> -3700.374899999999
$ws.Cells.Item(136, 14).Value = -15065.4  # N136: -16200 -> -15065.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2042.5834  # H3: 2155.6365 -> 2042.5834
$ws.Cells.Item(3, 9).Value = 1683.2727  # I3: 1771.7 -> 1683.2727
$ws.Cells.Item(3, 11).Value = 1683.2727  # K3: 1771.7 -> 1683.2727
$ws.Cells.Item(3, 13).Value = -1569.2727  # M3: -1657.7 -> -1569.2727
$ws.Cells.Item(94, 8).Value = 2921.6128  # H94: 3105.6333 -> 2921.6128
$ws.Cells.Item(94, 9).Value = 2728  # I94: 2886.2917 -> 2728
$ws.Cells.Item(94, 10).Value = 3585.4285  # J94: 3983 -> 3585.4285
$ws.Cells.Item(94, 11).Value = 2728  # K94: 2886.2917 -> 2728
$ws.Cells.Item(94, 12).Value = 3585.4285  # L94: 3983 -> 3585.4285
$ws.Cells.Item(94, 13).Value = -2277  # M94: -2435.2917 -> -2277
$ws.Cells.Item(94, 14).Value = -4487.4285  # N94: -4885 -> -4487.4285

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 80000  # H97: 55098.5 -> 80000
$ws.Cells.Item(97, 10).Value = 80000  # J97: 55098.5 -> 80000
$ws.Cells.Item(97, 12).Value = 80000  # L97: 55098.5 -> 80000
$ws.Cells.Item(97, 14).Value = -81982  # N97: -57080.5 -> -81982
$ws.Cells.Item(99, 8).Value = 37978.6  # H99: 40489.1 -> 37978.6
$ws.Cells.Item(99, 9).Value = 37978.6  # I99: 40489.1 -> 37978.6
$ws.Cells.Item(99, 11).Value = 37978.6  # K99: 40489.1 -> 37978.6
$ws.Cells.Item(99, 13).Value = -36480.6  # M99: -38991.1 -> -36480.6
$ws.Cells.Item(107, 8).Value = 1557.5  # H107: 1792.6522 -> 1557.5
$ws.Cells.Item(107, 9).Value = 1355.8823  # I107: 1615.5333 -> 1355.8823
$ws.Cells.Item(107, 10).Value = 1938.3334  # J107: 2124.75 -> 1938.3334
$ws.Cells.Item(107, 11).Value = 1355.8823  # K107: 1615.5333 -> 1355.8823
$ws.Cells.Item(107, 12).Value = 1938.3334  # L107: 2124.75 -> 1938.3334
$ws.Cells.Item(107, 13).Value = 564.1177  # M107: 304.4666999999999 -> 564.1177
$ws.Cells.Item(107, 14).Value = -5778.3334  # N107: -5964.75 -> -5778.3334
$ws.Cells.Item(122, 8).Value = 128106.625  # H122: 113975.78 -> 128106.625
$ws.Cells.Item(122, 9).Value = 204495.4  # I122: 128123.25 -> 204495.4
$ws.Cells.Item(122, 10).Value = 792  # J122: 796 -> 792
$ws.Cells.Item(122, 11).Value = 613486.2  # K122: 384369.75 -> 613486.2
$ws.Cells.Item(122, 12).Value = 2376  # L122: 2388 -> 2376
$ws.Cells.Item(122, 13).Value = -611036.2  # M122: -381919.75 -> -611036.2
$ws.Cells.Item(122, 14).Value = -7276  # N122: -7288 -> -7276
$ws.Cells.Item(126, 8).Value = 37978.6  # H126: 40489.1 -> 37978.6
$ws.Cells.Item(126, 9).Value = 37978.6  # I126: 40489.1 -> 37978.6
$ws.Cells.Item(126, 11).Value = 113935.8  # K126: 121467.3 -> 113935.8
$ws.Cells.Item(126, 13).Value = -111465.8  # M126: -118997.3 -> -111465.8
$ws.Cells.Item(134, 8).Value = 3063.1875  # H134: 2953.5293 -> 3063.1875
$ws.Cells.Item(134, 9).Value = 2934.0667  # I134: 2825.625 -> 2934.0667
$ws.Cells.Item(134, 11).Value = 8802.2001  # K134: 8476.875 -> 8802.2001
$ws.Cells.Item(134, 13).Value = -6267.2001  # M134: -5941.875 -> -6267.2001
$ws.Cells.Item(140, 8).Value = 62999.332  # H140: 55399.2 -> 62999.332
$ws.Cells.Item(140, 10).Value = 71999.5  # J140: 57999.25 -> 71999.5
$ws.Cells.Item(140, 12).Value = 71999.5  # L140: 57999.25 -> 71999.5
$ws.Cells.Item(140, 14).Value = -82359.5  # N140: -68359.25 -> -82359.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(45, 8).Value = 11873.5  # H45: 12142 -> 11873.5
$ws.Cells.Item(45, 10).Value = 11873.5  # J45: 12142 -> 11873.5
$ws.Cells.Item(45, 12).Value = 35620.5  # L45: 36426 -> 35620.5
$ws.Cells.Item(45, 14).Value = -36684.5  # N45: -37490 -> -36684.5
$ws.Cells.Item(70, 8).Value = 5498.8  # H70: 4664.8335 -> 5498.8
$ws.Cells.Item(70, 9).Value = 3997.5  # I70: 2245 -> 3997.5
$ws.Cells.Item(70, 10).Value = 6499.6665  # J70: 5874.75 -> 6499.6665
$ws.Cells.Item(70, 11).Value = 11992.5  # K70: 6735 -> 11992.5
$ws.Cells.Item(70, 12).Value = 19498.9995  # L70: 17624.25 -> 19498.9995
$ws.Cells.Item(70, 13).Value = -11677.5  # M70: -6420 -> -11677.5
$ws.Cells.Item(70, 14).Value = -20128.9995  # N70: -18254.25 -> -20128.9995
$ws.Cells.Item(73, 8).Value = 5498.8  # H73: 4664.8335 -> 5498.8
$ws.Cells.Item(73, 9).Value = 3997.5  # I73: 2245 -> 3997.5
$ws.Cells.Item(73, 10).Value = 6499.6665  # J73: 5874.75 -> 6499.6665
$ws.Cells.Item(73, 11).Value = 11992.5  # K73: 6735 -> 11992.5
$ws.Cells.Item(73, 12).Value = 19498.9995  # L73: 17624.25 -> 19498.9995
$ws.Cells.Item(73, 13).Value = -10900.5  # M73: -5643 -> -10900.5
$ws.Cells.Item(73, 14).Value = -21682.9995  # N73: -19808.25 -> -21682.9995
$ws.Cells.Item(75, 8).Value = 10580  # H75: 10750 -> 10580
$ws.Cells.Item(75, 10).Value = 10580  # J75: 10750 -> 10580
$ws.Cells.Item(75, 12).Value = 31740  # L75: 32250 -> 31740
$ws.Cells.Item(75, 14).Value = -33736  # N75: -34246 -> -33736
$ws.Cells.Item(78, 8).Value = 10580  # H78: 10750 -> 10580
$ws.Cells.Item(78, 10).Value = 10580  # J78: 10750 -> 10580
$ws.Cells.Item(78, 12).Value = 95220  # L78: 96750 -> 95220
$ws.Cells.Item(78, 14).Value = -105204  # N78: -106734 -> -105204
$ws.Cells.Item(107, 8).Value = 1185.5  # H107: 1223 -> 1185.5
$ws.Cells.Item(107, 9).Value = 899.3333  # I107: 1000 -> 899.3333
$ws.Cells.Item(107, 11).Value = 2697.9999  # K107: 3000 -> 2697.9999
$ws.Cells.Item(107, 13).Value = -777.9998999999998  # M107: -1080 -> -777.9998999999998
$ws.Cells.Item(114, 8).Value = 28573270  # H114: 16667978 -> 28573270
$ws.Cells.Item(114, 9).Value = 66667332  # I114: 33334006 -> 66667332
$ws.Cells.Item(114, 10).Value = 2722.25  # J114: 1949.5 -> 2722.25
$ws.Cells.Item(114, 11).Value = 200001996  # K114: 100002018 -> 200001996
$ws.Cells.Item(114, 12).Value = 8166.75  # L114: 5848.5 -> 8166.75
$ws.Cells.Item(114, 13).Value = -199998742  # M114: -99998764 -> -199998742
$ws.Cells.Item(114, 14).Value = -14674.75  # N114: -12356.5 -> -14674.75
$ws.Cells.Item(122, 8).Value = 527.8182  # H122: 517.7273 -> 527.8182
$ws.Cells.Item(122, 10).Value = 494.5  # J122: 466.75 -> 494.5
$ws.Cells.Item(122, 12).Value = 4450.5  # L122: 4200.75 -> 4450.5
$ws.Cells.Item(122, 14).Value = -9350.5  # N122: -9100.75 -> -9350.5
$ws.Cells.Item(131, 8).Value = 10918.454  # H131: 11535.3 -> 10918.454
$ws.Cells.Item(131, 9).Value = 12650.333  # I131: 11535.3 -> 12650.333
$ws.Cells.Item(131, 10).Value = 3125  # J131: 0 -> 3125
$ws.Cells.Item(131, 11).Value = 37950.999  # K131: 34605.89999999999 -> 37950.999
$ws.Cells.Item(131, 12).Value = 9375  # L131: 0 -> 9375
$ws.Cells.Item(131, 13).Value = -32910.999  # M131: -29565.89999999999 -> -32910.999
$ws.Cells.Item(131, 14).Value = -19455  # N131: None -> -19455
$ws.Cells.Item(132, 8).Value = 2194.4167  # H132: 2164.72 -> 2194.4167
$ws.Cells.Item(132, 10).Value = 2320.389  # J132: 2274.6843 -> 2320.389
$ws.Cells.Item(132, 12).Value = 20883.501  # L132: 20472.1587 -> 20883.501
$ws.Cells.Item(132, 14).Value = -25943.501  # N132: -25532.1587 -> -25943.501

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2965.0833  # H80: 3059.8635 -> 2965.0833
$ws.Cells.Item(80, 9).Value = 2451.4707  # I80: 2559.3572 -> 2451.4707
$ws.Cells.Item(80, 10).Value = 4212.4287  # J80: 3935.75 -> 4212.4287
$ws.Cells.Item(80, 11).Value = 2451.4707  # K80: 2559.3572 -> 2451.4707
$ws.Cells.Item(80, 12).Value = 4212.4287  # L80: 3935.75 -> 4212.4287
$ws.Cells.Item(80, 13).Value = -1453.4707  # M80: -1561.3572 -> -1453.4707
$ws.Cells.Item(80, 14).Value = -6208.4287  # N80: -5931.75 -> -6208.4287
$ws.Cells.Item(83, 8).Value = 2965.0833  # H83: 3059.8635 -> 2965.0833
$ws.Cells.Item(83, 9).Value = 2451.4707  # I83: 2559.3572 -> 2451.4707
$ws.Cells.Item(83, 10).Value = 4212.4287  # J83: 3935.75 -> 4212.4287
$ws.Cells.Item(83, 11).Value = 12257.3535  # K83: 12796.786 -> 12257.3535
$ws.Cells.Item(83, 12).Value = 21062.1435  # L83: 19678.75 -> 21062.1435
$ws.Cells.Item(83, 13).Value = -7265.353499999999  # M83: -7804.786 -> -7265.353499999999
$ws.Cells.Item(83, 14).Value = -31046.1435  # N83: -29662.75 -> -31046.1435
$ws.Cells.Item(122, 8).Value = 1087.2609  # H122: 1201.15 -> 1087.2609
$ws.Cells.Item(122, 9).Value = 1130.5  # I122: 1245.7778 -> 1130.5
$ws.Cells.Item(122, 10).Value = 799  # J122: 799.5 -> 799
$ws.Cells.Item(122, 11).Value = 3391.5  # K122: 3737.3334 -> 3391.5
$ws.Cells.Item(122, 12).Value = 2397  # L122: 2398.5 -> 2397
$ws.Cells.Item(122, 13).Value = -941.5  # M122: -1287.3334 -> -941.5
$ws.Cells.Item(122, 14).Value = -7297  # N122: -7298.5 -> -7297

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 13078.849  # H7: 13696.742 -> 13078.849
$ws.Cells.Item(7, 9).Value = 41000.285  # I7: 55999.8 -> 41000.285
$ws.Cells.Item(7, 11).Value = 41000.285  # K7: 55999.8 -> 41000.285
$ws.Cells.Item(7, 13).Value = -40888.285  # M7: -55887.8 -> -40888.285
$ws.Cells.Item(61, 8).Value = 118304.36  # H61: 137786.75 -> 118304.36
$ws.Cells.Item(61, 9).Value = 106077.336  # I61: 127010.8 -> 106077.336
$ws.Cells.Item(61, 11).Value = 106077.336  # K61: 127010.8 -> 106077.336
$ws.Cells.Item(61, 13).Value = -105875.336  # M61: -126808.8 -> -105875.336
$ws.Cells.Item(82, 8).Value = 3892.375  # H82: 4065.2 -> 3892.375
$ws.Cells.Item(82, 9).Value = 2954.818  # I82: 3120.3 -> 2954.818
$ws.Cells.Item(82, 11).Value = 2954.818  # K82: 3120.3 -> 2954.818
$ws.Cells.Item(82, 13).Value = -2593.818  # M82: -2759.3 -> -2593.818
$ws.Cells.Item(85, 8).Value = 3892.375  # H85: 4065.2 -> 3892.375
$ws.Cells.Item(85, 9).Value = 2954.818  # I85: 3120.3 -> 2954.818
$ws.Cells.Item(85, 11).Value = 2954.818  # K85: 3120.3 -> 2954.818
$ws.Cells.Item(85, 13).Value = -1706.818  # M85: -1872.3 -> -1706.818
$ws.Cells.Item(113, 8).Value = 118304.36  # H113: 137786.75 -> 118304.36
$ws.Cells.Item(113, 9).Value = 106077.336  # I113: 127010.8 -> 106077.336
$ws.Cells.Item(113, 11).Value = 106077.336  # K113: 127010.8 -> 106077.336
$ws.Cells.Item(113, 13).Value = -103907.336  # M113: -124840.8 -> -103907.336
$ws.Cells.Item(122, 8).Value = 5263.4653  # H122: 5263.4883 -> 5263.4653
$ws.Cells.Item(122, 9).Value = 5627  # I122: 5627.032 -> 5627
$ws.Cells.Item(122, 11).Value = 16881  # K122: 16881.096 -> 16881
$ws.Cells.Item(122, 13).Value = -14431  # M122: -14431.096 -> -14431
$ws.Cells.Item(126, 8).Value = 13078.849  # H126: 13696.742 -> 13078.849
$ws.Cells.Item(126, 9).Value = 41000.285  # I126: 55999.8 -> 41000.285
$ws.Cells.Item(126, 11).Value = 123000.855  # K126: 167999.4 -> 123000.855
$ws.Cells.Item(126, 13).Value = -120530.855  # M126: -165529.4 -> -120530.855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 853.6087  # H126: 767.8077 -> 853.6087
$ws.Cells.Item(126, 9).Value = 687.7222  # I126: 605.1905 -> 687.7222
$ws.Cells.Item(126, 11).Value = 2063.1666  # K126: 1815.5715 -> 2063.1666
$ws.Cells.Item(126, 13).Value = 406.8334  # M126: 654.4285 -> 406.8334
$ws.Cells.Item(132, 8).Value = 8019.3335  # H132: 7277.95 -> 8019.3335
$ws.Cells.Item(132, 9).Value = 11219.182  # I132: 9115 -> 11219.182
$ws.Cells.Item(132, 10).Value = 2991  # J132: 2991.5 -> 2991
$ws.Cells.Item(132, 11).Value = 33657.546  # K132: 27345 -> 33657.546
$ws.Cells.Item(132, 12).Value = 8973  # L132: 8974.5 -> 8973
$ws.Cells.Item(132, 13).Value = -31127.546  # M132: -24815 -> -31127.546
$ws.Cells.Item(132, 14).Value = -14033  # N132: -14034.5 -> -14033
